# Revert "adding additional models"
# Restores "Combined Results" and "Neural Network" sheets to their
# pre-feature state (fewer rows, F1-score summary columns on Combined
# Results, and an "Unscaled" scaler block restored on Neural Network).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: Neural Network  (edited first so the final active tab/selection
# ends up on "Combined Results", matching the target workbook view)
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Neural Network")

$ws5.Range("A1:D14").ClearContents()

$ws5.Cells.Item(1, 1).Value = "Dataset"
$ws5.Cells.Item(1, 2).Value = "Scaler"
$ws5.Cells.Item(1, 3).Value = "Features"
$ws5.Cells.Item(1, 4).Value = "Neural Network"

$headerCell5 = $ws5.Cells.Item(1, 1)
$headerCell5.Font.Bold = $true
$headerCell5.Borders.Item(9).LineStyle = 1
$headerCell5.Copy()
$ws5.Range("B1:G1").PasteSpecial(-4122)

$ws5.Cells.Item(2, 1).Value = "Actual"
$ws5.Cells.Item(2, 2).Value = "Min-Max"
$ws5.Cells.Item(2, 3).Value = "All"

$ws5.Cells.Item(3, 2).Value = "Min-Max"
$ws5.Cells.Item(3, 3).Value = "Med Only"

$ws5.Cells.Item(4, 2).Value = "Standard"
$ws5.Cells.Item(4, 3).Value = "All"

$ws5.Cells.Item(5, 2).Value = "Standard"
$ws5.Cells.Item(5, 3).Value = "Med Only"

$ws5.Cells.Item(6, 2).Value = "Unscaled"
$ws5.Cells.Item(6, 3).Value = "All"

$ws5.Cells.Item(7, 2).Value = "Unscaled"
$ws5.Cells.Item(7, 3).Value = "Med Only"

$ws5.Cells.Item(8, 1).Value = "RandomOverSampled"
$ws5.Cells.Item(8, 2).Value = "Min-Max"
$ws5.Cells.Item(8, 3).Value = "All"

$ws5.Cells.Item(9, 2).Value = "Min-Max"
$ws5.Cells.Item(9, 3).Value = "Med Only"

$ws5.Cells.Item(10, 2).Value = "Standard"
$ws5.Cells.Item(10, 3).Value = "All"

$ws5.Cells.Item(11, 2).Value = "Standard"
$ws5.Cells.Item(11, 3).Value = "Med Only"

$ws5.Cells.Item(12, 1).Value = "SMOTEENN"
$ws5.Cells.Item(12, 2).Value = "Min-Max"
$ws5.Cells.Item(12, 3).Value = "All"

$ws5.Cells.Item(13, 2).Value = "Min-Max"
$ws5.Cells.Item(13, 3).Value = "Med Only"

$ws5.Cells.Item(14, 2).Value = "Standard"
$ws5.Cells.Item(14, 3).Value = "All"

$ws5.Cells.Item(15, 2).Value = "Standard"
$ws5.Cells.Item(15, 3).Value = "Med Only"

$ws5.PageSetup.Orientation = 1
$ws5.Range("E7").Select()

# ---------------------------------------------------------------------
# Sheet: Combined Results
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Combined Results")

# Wipe the old layout (title cell in row 1 + data starting row 2).
$ws1.Range("A1:G14").ClearContents()

# New header row (bold, bottom-border).
$ws1.Cells.Item(1, 1).Value = "Dataset"
$ws1.Cells.Item(1, 2).Value = "Scaler"
$ws1.Cells.Item(1, 3).Value = "Features"
$ws1.Cells.Item(1, 4).Value = "Neural Network"
$ws1.Cells.Item(1, 5).Value = "Random Forest"
$ws1.Cells.Item(1, 6).Value = "Logistic Regression"
$ws1.Cells.Item(1, 7).Value = "Decision Tree"

$headerCell1 = $ws1.Cells.Item(1, 1)
$headerCell1.Font.Bold = $true
$headerCell1.Borders.Item(9).LineStyle = 1
$headerCell1.Copy()
$ws1.Range("B1:G1").PasteSpecial(-4122)

# Data rows.
$ws1.Cells.Item(2, 1).Value = "Actual"
$ws1.Cells.Item(2, 2).Value = "Standard"
$ws1.Cells.Item(2, 3).Value = "All"
$ws1.Cells.Item(2, 5).Value = "F1 = 0.80 / 0.82"
$ws1.Cells.Item(2, 6).Value = "F1 = 0.84 / 0.85"

$ws1.Cells.Item(3, 2).Value = "Standard"
$ws1.Cells.Item(3, 3).Value = "Med Only"
$ws1.Cells.Item(3, 5).Value = "F1 = 0.82 /0.82"
$ws1.Cells.Item(3, 6).Value = "F1 = 0.84 / 0.85"

$ws1.Cells.Item(4, 2).Value = "Unscaled"
$ws1.Cells.Item(4, 3).Value = "All"
$ws1.Cells.Item(4, 5).Value = "F1 = 0.80 / 0.82"
$ws1.Cells.Item(4, 6).Value = "F1 = 0.84 / 0.85"

$ws1.Cells.Item(5, 2).Value = "Unscaled"
$ws1.Cells.Item(5, 3).Value = "Med Only"
$ws1.Cells.Item(5, 5).Value = "F1 = 0.83 / 0.83"
$ws1.Cells.Item(5, 6).Value = "F1 = 0.84 / 0.85"

$ws1.Cells.Item(6, 1).Value = "SMOTEENN"
$ws1.Cells.Item(6, 2).Value = "Standard"
$ws1.Cells.Item(6, 3).Value = "All"
$ws1.Cells.Item(6, 5).Value = "F1 = 0.94 / 0.24"
$ws1.Cells.Item(6, 6).Value = "F1 = 0.82/ 0.21"

$ws1.Cells.Item(7, 2).Value = "Standard"
$ws1.Cells.Item(7, 3).Value = "Med Only"
$ws1.Cells.Item(7, 5).Value = "F1 = 0.93 / 0.24"
$ws1.Cells.Item(7, 6).Value = "F1 = 0.82/ 0.21"

# Final selection / active sheet: "Combined Results", cell F6 — this is
# the last thing the script touches so it becomes the active tab.
$ws1.Range("F6").Select()
